$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.209.35"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "2.248.00"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "246.97"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.19%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.630"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "74.52"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -5.90%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.619"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -5.29%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "42.09"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("E11").Value = "  -2.92%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("E13").Value = "  -1.93%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "14.58"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "2.249.07"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "42.128.67"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "0.0₃0990"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  +0.39%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.23"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.08%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "231.53"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.92"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +37.85%  "
$ws.Range("E24").Value = "  +0.07%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "11.42"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  -5.03%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.31"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  -0.24%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "168.89"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  -1.06%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0829"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("E32").Value = "  +0.50%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "31.33"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("E34").Value = "  -1.88%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.25"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +9.17%  "
$ws.Range("E36").Value = "  -1.58%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0315"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.49%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "13.91"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.51%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.18"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.64%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "62.57"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.205"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "105.72"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -5.68%  "
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.17"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.12"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.94%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.17"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -10.03%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.69"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.20%  "
